# Scheduled market-data refresh: push freshly-pulled price/profit figures
# into each job sheet (Leve item market data + derived NQ/HQ profit calcs).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2867.1667
$ws.Range("I8").Value = 2867.1667
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 8601.500100000001
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -8462.500100000001
$ws.Range("N8").ClearContents()

$ws.Range("H86").Value = 9478.333000000001
$ws.Range("I86").Value = 9521.862999999999
$ws.Range("K86").Value = 9521.862999999999
$ws.Range("M86").Value = -8398.862999999999

$ws.Range("H89").Value = 9478.333000000001
$ws.Range("I89").Value = 9521.862999999999
$ws.Range("K89").Value = 47609.315
$ws.Range("M89").Value = -41993.315

$ws.Range("H132").Value = 2247.4792
$ws.Range("I132").Value = 1805.317
$ws.Range("J132").Value = 4837.2856
$ws.Range("K132").Value = 5415.951
$ws.Range("L132").Value = 14511.8568
$ws.Range("M132").Value = -2885.951
$ws.Range("N132").Value = -19571.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7972.8965
$ws.Range("I132").Value = 3663.2942
$ws.Range("K132").Value = 10989.8826
$ws.Range("M132").Value = -8459.882599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1441.6
$ws.Range("I107").Value = 1458.5555
$ws.Range("K107").Value = 1458.5555
$ws.Range("M107").Value = 461.4445000000001

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2227.75
$ws.Range("I16").Value = 2303.6667
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 2303.6667
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -2016.6667
$ws.Range("N16").Value = -2574

$ws.Range("H31").Value = 1071738.8
$ws.Range("I31").Value = 32399.4
$ws.Range("J31").Value = 1377426.8
$ws.Range("K31").Value = 32399.4
$ws.Range("L31").Value = 1377426.8
$ws.Range("M31").Value = -32104.4
$ws.Range("N31").Value = -1378016.8

$ws.Range("H34").Value = 1071738.8
$ws.Range("I34").Value = 32399.4
$ws.Range("J34").Value = 1377426.8
$ws.Range("K34").Value = 32399.4
$ws.Range("L34").Value = 1377426.8
$ws.Range("M34").Value = -32197.4
$ws.Range("N34").Value = -1377830.8

$ws.Range("H62").Value = 2653.7273
$ws.Range("I62").Value = 2520.8
$ws.Range("J62").Value = 3983
$ws.Range("K62").Value = 2520.8
$ws.Range("L62").Value = 3983
$ws.Range("M62").Value = -1896.8
$ws.Range("N62").Value = -5231

$ws.Range("H65").Value = 2653.7273
$ws.Range("I65").Value = 2520.8
$ws.Range("J65").Value = 3983
$ws.Range("K65").Value = 12604
$ws.Range("L65").Value = 19915
$ws.Range("M65").Value = -9484
$ws.Range("N65").Value = -26155

$ws.Range("H105").Value = 2427.182
$ws.Range("I105").Value = 2257
$ws.Range("J105").Value = 2491
$ws.Range("K105").Value = 2257
$ws.Range("L105").Value = 2491
$ws.Range("M105").Value = -510
$ws.Range("N105").Value = -5985

$ws.Range("H111").Value = 59996.332
$ws.Range("J111").Value = 59996.332
$ws.Range("L111").Value = 59996.332
$ws.Range("N111").Value = -68176.33199999999

$ws.Range("H113").Value = 2227.75
$ws.Range("I113").Value = 2303.6667
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2303.6667
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -133.6667000000002
$ws.Range("N113").Value = -6340

$ws.Range("H118").Value = 115000
$ws.Range("J118").Value = 115000
$ws.Range("L118").Value = 115000
$ws.Range("N118").Value = -118314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1273.909
$ws.Range("I113").Value = 845
$ws.Range("J113").Value = 1788.6
$ws.Range("K113").Value = 2535
$ws.Range("L113").Value = 5365.799999999999
$ws.Range("M113").Value = -365
$ws.Range("N113").Value = -9705.799999999999

$ws.Range("H129").Value = 854.7
$ws.Range("I129").Value = 567.63635
$ws.Range("J129").Value = 1205.5555
$ws.Range("K129").Value = 1702.90905
$ws.Range("L129").Value = 3616.6665
$ws.Range("M129").Value = 3297.09095
$ws.Range("N129").Value = -13616.6665

$ws.Range("H134").Value = 3726.3809
$ws.Range("I134").Value = 2391
$ws.Range("K134").Value = 7173
$ws.Range("M134").Value = -2103

$ws.Range("H140").Value = 275231.72
$ws.Range("I140").Value = 302339.9
$ws.Range("K140").Value = 907019.7000000001
$ws.Range("M140").Value = -901839.7000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H109").Value = 46774.5
$ws.Range("J109").Value = 46774.5
$ws.Range("L109").Value = 46774.5
$ws.Range("N109").Value = -48854.5

$ws.Range("H113").Value = 3404.6667
$ws.Range("I113").Value = 3240
$ws.Range("J113").Value = 4475
$ws.Range("K113").Value = 3240
$ws.Range("L113").Value = 4475
$ws.Range("M113").Value = -1070
$ws.Range("N113").Value = -8815

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 44378.88
$ws.Range("I7").Value = 3897.1667
$ws.Range("J7").Value = 148474.72
$ws.Range("K7").Value = 3897.1667
$ws.Range("L7").Value = 148474.72
$ws.Range("M7").Value = -3785.1667
$ws.Range("N7").Value = -148698.72

$ws.Range("H22").Value = 2462.5186
$ws.Range("I22").Value = 2474.4
$ws.Range("J22").Value = 2428.5715
$ws.Range("K22").Value = 2474.4
$ws.Range("L22").Value = 2428.5715
$ws.Range("M22").Value = -2179.4
$ws.Range("N22").Value = -3018.5715

$ws.Range("H27").Value = 2462.5186
$ws.Range("I27").Value = 2474.4
$ws.Range("J27").Value = 2428.5715
$ws.Range("K27").Value = 2474.4
$ws.Range("L27").Value = 2428.5715
$ws.Range("M27").Value = -2367.4
$ws.Range("N27").Value = -2642.5715

$ws.Range("H68").Value = 1999.5
$ws.Range("I68").Value = 1999.5
$ws.Range("K68").Value = 1999.5
$ws.Range("M68").Value = -1250.5

$ws.Range("H71").Value = 1999.5
$ws.Range("I71").Value = 1999.5
$ws.Range("K71").Value = 9997.5
$ws.Range("M71").Value = -6253.5

$ws.Range("H101").Value = 68872.39999999999
$ws.Range("J101").Value = 68872.39999999999
$ws.Range("L101").Value = 68872.39999999999
$ws.Range("N101").Value = -75362.39999999999

$ws.Range("H126").Value = 44378.88
$ws.Range("I126").Value = 3897.1667
$ws.Range("J126").Value = 148474.72
$ws.Range("K126").Value = 11691.5001
$ws.Range("L126").Value = 445424.16
$ws.Range("M126").Value = -9221.500100000001
$ws.Range("N126").Value = -450364.16

$ws.Range("H136").Value = 38091.547
$ws.Range("I136").Value = 1905.7826
$ws.Range("J136").Value = 142125.62
$ws.Range("K136").Value = 5717.3478
$ws.Range("L136").Value = 426376.86
$ws.Range("M136").Value = -3167.3478
$ws.Range("N136").Value = -431476.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H136").Value = 8037.5835
$ws.Range("I136").Value = 1095.32
$ws.Range("K136").Value = 3285.96
$ws.Range("M136").Value = -735.96
